$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B12").Value = 140809
$ws.Range("B15").Value = 47
$ws.Range("C15").Value = "Cluj-Cluj"
$ws.Range("B16").Value = 30
$ws.Range("C16").Value = "Acasa-Birou"
$ws.Range("D16").Value = " "
$ws.Range("B17").Value = 30
$ws.Range("C17").Value = "Acasa-Birou"
$ws.Range("D17").Value = " "
$ws.Range("B19").Value = 121
$ws.Range("C19").Value = "Cluj-Turda"
$ws.Range("B22").Value = 30
$ws.Range("C22").Value = "Acasa-Birou"
$ws.Range("D22").Value = " "
$ws.Range("B23").Value = 101
$ws.Range("C23").Value = "Cluj-Dej"
$ws.Range("B24").Value = 47
$ws.Range("C24").Value = "Cluj-Cluj"
$ws.Range("B25").Value = 152
$ws.Range("C25").Value = "Cluj-Cmp. Turzii"
$ws.Range("D25").Value = "Interes Serviciu"
$ws.Range("B26").Value = 257
$ws.Range("C26").Value = "Cluj-Bistrita"
$ws.Range("B29").Value = 30
$ws.Range("C29").Value = "Acasa-Birou"
$ws.Range("D29").Value = " "
$ws.Range("B30").Value = 121
$ws.Range("C30").Value = "Cluj-Turda"
$ws.Range("D30").Value = "Interes Serviciu"
$ws.Range("B33").Value = 30
$ws.Range("C33").Value = "Acasa-Birou"
$ws.Range("D33").Value = " "
$ws.Range("B36").Value = 85
$ws.Range("C36").Value = "Cluj-Apahida"
$ws.Range("D36").Value = "Interes Serviciu"
$ws.Range("B38").Value = 156
$ws.Range("C38").Value = "Cluj-Zalau"
$ws.Range("B39").Value = 257
$ws.Range("C39").Value = "Cluj-Bistrita"
$ws.Range("D39").Value = "Interes Serviciu"
$ws.Range("B40").Value = 121
$ws.Range("C40").Value = "Cluj-Turda"
$ws.Range("B43").Value = 30
$ws.Range("C43").Value = "Acasa-Birou"
$ws.Range("D43").Value = " "
$ws.Range("B44").Value = 85
$ws.Range("C44").Value = "Cluj-Apahida"
$ws.Range("B45").Value = 1820
$ws.Range("B46").Value = 142629
